# Trade #23 closed at 2026-02-17 13:18:36 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" aggregate numbers to reflect
# the newly closed trade, and appends the new trade row (row 24) to both
# the "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate metrics
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.17   # Current Capital
$wsSummary.Range("B4").Value = -0.84     # Total P&L $
$wsSummary.Range("B5").Value = -0.73     # Total P&L %
$wsSummary.Range("B6").Value = 23        # Total Trades
$wsSummary.Range("B8").Value = 15        # Losing Trades
$wsSummary.Range("B9").Value = 30.43     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: refresh the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.17      # Capital
$wsStatus.Range("D4").Value = 23         # Trades
$wsStatus.Range("E4").Value = -0.84      # P&L $
$wsStatus.Range("F4").Value = -0.83      # P&L %
$wsStatus.Range("G4").Value = 30.43      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade (#23) to a worksheet's trade log, as row 24.
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 24

    $ws.Cells.Item($row, 1).Value = 23

    # Date column looks like a date ("2026-02-17"); force Text formatting
    # first so Excel stores the literal string instead of auto-converting
    # it to a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "13:18:29"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.36
    $ws.Cells.Item($row, 7).Value = 0.29
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -19.4444
    $ws.Cells.Item($row, 10).Value = -0.07000000000000001
    $ws.Cells.Item($row, 11).Value = 99.17
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
